# Apply the CDA Logical Model update for ManufacturedProduct (ST.r2b release)
#  - bump Version / Date metadata
#  - add a new "Jurisdiction" property row to the Metadata sheet
#  - document the II-1 constraint on ManufacturedProduct.typeId in the Elements sheet

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Date bump
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" row right after "Contact" (row 10), before "Description" (row 11)
$meta.Rows.Item(11).Insert()

# Copy formatting from the row below (still a normal data row) onto the newly inserted row
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# Document the II-1 invariant on ManufacturedProduct.typeId (row 5, column AJ = Constraint(s))
$elements.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}`n"
